# Fill in the "General Report" and "Blockers" standup bullet lists.
$d = $word.ActiveDocument

# --- "General Report" section (bulleted list, numId=3, ilvl=0) ---
# The first bullet under "General Report" is currently an empty list
# paragraph; give it Eric's update, then append two siblings (same list
# formatting) for Sam and Luke.
$general = $d.Paragraphs.Item(3)
$general.Range.Text = "Eric is working on mastering Python Classes"

$general.Range.InsertParagraphAfter()
$sam1 = $d.Paragraphs.Item(4)
$sam1.Range.Text = "Sam is working on understanding Python Functions"

$sam1.Range.InsertParagraphAfter()
$luke = $d.Paragraphs.Item(5)
$luke.Range.Text = "Luke is studying how to work with String objects"

# --- "Blockers" section (bulleted list, numId=2) ---
# The lone empty bullet under "Blockers" is now paragraph 7 (two new
# paragraphs were inserted above). Give it Eric's blocker, then add a
# second, more deeply indented bullet for Sam's follow-up.
$blocker = $d.Paragraphs.Item(7)
$blocker.Range.Text = "Eric is struggling to understand what static methods are used for in classes"

$blocker.Range.InsertParagraphAfter()
$followUp = $d.Paragraphs.Item(8)
$followUp.Range.ListFormat.ListLevelNumber = 2
$followUp.Range.Text = "Sam helped him understand they are meant to be used as utility methods"
